# edit.ps1
# Update "Market Health Data" workbook per commit: Update data: 2025-11-05 12:41
#
# Changes:
# 1. Metadata!A2 timestamp bumped from 12:40 PM to 12:41 PM
# 2. "Industry Analysis" sheet column F ("1 Year") values refreshed for rows 2-76
# 3. "Stock List" sheet rows 2-76 shifted up by one (new top-of-list entry removed,
#    a new entry "TRAVELFOOD" appended at the bottom) for columns B, C (stock code),
#    D (Price), E (% Change) and H (Market Cap)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Metadata sheet - bump the "Last Updated" timestamp
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "05 Nov 2025, 12:41 PM"

# ---------------------------------------------------------------------------
# 2) Industry Analysis sheet - refresh "1 Year" (column F) values, rows 2-76
# ---------------------------------------------------------------------------
$industry = $wb.Worksheets.Item("Industry Analysis")

$industryData = @"
2|18.476
3|-7.7404
4|30.7972
5|-50.2266
6|61.9649
7|-9.1713
8|-3.556
9|38.3509
10|-6.2497
11|52.6723
12|-6.932
13|17.5662
14|-35.5106
15|0.6286
16|-3.1514
17|-20.6354
18|-0.0175
19|-26.9255
20|44.703
21|10.0506
22|84.6016
23|-54.4868
24|-12.8122
25|-9.182700000000001
26|5.9529
27|-33.2998
28|-20.4441
29|-17.1514
30|24.527
31|57.6193
32|-1.527
33|-5.2378
34|27.4054
35|6.7961
36|-5.6683
37|1.4178
38|-22.4272
39|12.3741
40|-5.138
41|-0.1825
42|23.2483
43|14.456
44|-11.1739
45|27.112
46|-5.6252
47|-36.5148
48|-27.8397
49|-25.4424
50|-49.1173
51|-51.065
52|-35.4517
53|-11.9879
54|-3.0992
55|-15.3441
56|-25.937
57|-29.1486
58|-6.4093
59|-23.3046
60|-11.2657
61|-9.777699999999999
62|-16.0561
63|-9.932499999999999
64|51.8767
65|-43.5191
66|13.7315
67|12.6111
68|31.7532
69|-19.9577
70|-12.9642
71|13.2432
72|2.8232
73|-9.179
74|-14.2931
75|28.3699
76|45.5868
"@

foreach ($line in ($industryData -split "`n")) {
    $line = $line.Trim()
    if ($line.Length -eq 0) { continue }
    $parts = $line -split '\|'
    $row = [int]$parts[0]
    $val = [double]$parts[1]
    $industry.Cells.Item($row, 6).Value = $val
}

# ---------------------------------------------------------------------------
# 3) Stock List sheet - refresh rows 2-76 (list shifted up by one entry, with
#    a new "TRAVELFOOD" row appended at the bottom)
# ---------------------------------------------------------------------------
$stocklist = $wb.Worksheets.Item("Stock List")

$stockData = @"
2|NIFTYCASE|10.19|-0.5854|0
3|MOMENTUM30|31.54|-0.6614|0
4|CANHLIFE|118.46|0.6286|11253.7
5|FLEXIADD|10.64|-1.0233|0
6|MOENERGY|36.3|-0.6568000000000001|0
7|MONIFTY100|26.49|0.3409|0
8|RUBICON|652.65|-0.1453|10752.4289
9|CRAMC|317.2|2.3226|6325.5208
10|LGEINDIA|1633.4|-0.946|110870.6825
11|TATACAP|329.3|0.1521|139783.5374
12|ELIQUID|1004.85|0.0408|0
13|WEWORK|632.15|-2.4008|8472.2803
14|GROWWRLTY|10.8|-0.4608|0
15|ADVANCE|130.05|-5.2666|836.0358
16|OMFREIGHT|88.90000000000001|-0.5926|299.3747
17|GLOTTIS|72.73999999999999|-0.8587|672.1394
18|FABTECH|237.72|0.4734|1056.6843
19|PACEDIGITK|218.85|0.1327|4723.9063
20|JAINREC|377.25|1.2208|13018.3623
21|EPACKPEB|301.45|1.979|3028.1254
22|BMWVENTLTD|69.25|0|600.5014
23|STYL|372.4|-0.8388|6025.649
24|JARO|621.5|-1.4821|1377.0134
25|SOLARWORLD|309.1|-0.6269|2679.0517
26|ARSSBL|537.3|4.7266|3370.2277
27|GANESHCP|274.4|-2.7984|1108.9312
28|ATLANTAELE|1003.05|-1.7436|7713.116
29|GKENERGY|213.85|-0.7933|4337.2472
30|SAATVIKGL|528.2|-1.3079|6713.6863
31|IVALUE|281.45|-0.3364|1506.8799
32|VMSTMT|70.03|-0.9056|347.5674
33|EUROPRATIK|321.75|0.8147|3288.285
34|SHRINGARMS|229.31|-1.2616|2211.284
35|DEVX|44.53|-0.3803|401.605
36|URBANCO|148.9|-2.0459|21380.5798
37|SML100CASE|10.36|-0.7663|0
38|AONEGOLD|11.28|-0.2653|0
39|ELM250|16.72|0.1797|0
40|AMANTA|122.52|1.407|475.7372
41|CPEDU|315.9|1.8539|574.7148999999999
42|AHCL|139.27|3.1706|740.2409
43|STLNETWORK|26.59|-0.412|1297.3822
44|VIKRAN|98.05|-1.783|2528.8166
45|MANUFGBEES|151.77|-1.011|0
46|MEIL|461.15|-0.7319|1274.1632
47|GROWWNXT50|70.29000000000001|-0.4109|0
48|SHREEJISPG|270.05|-0.7899|4399.6074
49|GEMAROMA|219.52|-0.876|1146.7097
50|PATELRMART|219.31|-1.0646|732.5069999999999
51|VIKRAMSOLR|322|-1.5892|11647.2884
52|LTGILTCASE|29.67|0.2365|0
53|REGAAL|89.13|-0.8675|915.5742
54|BLUESTONE|711.95|0.1266|10773.2539
55|MOSILVER|145.9|-1.5054|0
56|ALLTIME|308.75|2.66|2022.5526
57|JSWCEMENT|134.98|-0.4793|18402.6999
58|SBILIQETF|1012.94|0.0296|0
59|HILINFRA|77.23|-0.3998|0
60|GROWWPOWER|10.28|-0.9634|0
61|LOTUSDEV|177.82|0.3669|8690.485000000001
62|MBEL|450.2|-0.7714|2572.8126
63|LAXMIINDIA|145.62|-1.1942|761.1248000000001
64|CPPLUS|1322.1|-0.264|15497.9053
65|SHANTIGOLD|241.57|-1.6409|1741.6231
66|MOGOLD|119.65|-0.5403|0
67|BRIGHOTEL|82.39|-0.9855|3129.5229
68|INDIQUBE|212.64|-0.7561|4465.6847
69|EBGNG|346.65|3.2311|3952.2092
70|LIQGRWBEES|1014.74|0.0246|0
71|CHEMBONDCH|153.35|-1.6987|412.459
72|GROWWNIFTY|10.29|-0.3872|0
73|ANTHEM|702.25|-0.1209|39439.0658
74|QUALITY30|21.05|-0.8945|0
75|SMARTWORKS|606.65|2.0867|6931.2448
76|TRAVELFOOD|1316.3|0.1141|17332.9705
"@

foreach ($line in ($stockData -split "`n")) {
    $line = $line.Trim()
    if ($line.Length -eq 0) { continue }
    $parts = $line -split '\|'
    $row = [int]$parts[0]
    $code = $parts[1]
    $price = [double]$parts[2]
    $pctChange = [double]$parts[3]
    $marketCap = [double]$parts[4]

    $stocklist.Cells.Item($row, 2).Value = $code
    $stocklist.Cells.Item($row, 3).Value = $code
    $stocklist.Cells.Item($row, 4).Value = $price
    $stocklist.Cells.Item($row, 5).Value = $pctChange
    $stocklist.Cells.Item($row, 8).Value = $marketCap
}
